$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column D: BER values for phase offset pi/2 (similar to Divsalar & Simon)
$ws.Range("D3").Formula = "=0.0034"
$ws.Range("D4").Formula = "=0.0009999"
$ws.Range("D5").Formula = "=0.00019"
$ws.Range("D6").Formula = "=0.00001"

# Reposition the chart (moved up/left on the sheet)
$co = $ws.ChartObjects().Item(1)
$co.Left = 269.40385826771654
$co.Top = 9.346141732283465
$co.Width = 444.4230708661417
$co.Height = 216.69236220472442

$ws.Range("D7").Select()
